$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: mark "Section 3 & 4" as complete (progress flag 0 -> 1) ---
$ws.Range("D26").Value = 1

# --- Row 27: mark "Section 6" complete, push its end date out a day (5 -> 6 days) ---
$ws.Range("D27").Value = 1
$ws.Range("F27").Formula = "=E25+5"

# --- Row 28: mark "Overall review/revision" complete, push its end date out a day (1 -> 2 days) ---
$ws.Range("D28").Value = 1
$ws.Range("F28").Formula = "=E25+5"

# --- Row 30: "get feedback from tutor" now starts 12 days after F28 instead of 15, and is not done ---
$ws.Range("D30").Value = 0
$ws.Range("E30").Formula = "=F28+12"

# --- Row 31: rename task, mark complete, and move it right after F28 ---
$ws.Range("B31").Value = "cloning summary"
$ws.Range("D31").Value = 1
$ws.Range("E31").Formula = "=F28+1"

# --- Insert a brand-new task row at 32: "add source codes" ---
$ws.Rows(32).Insert()

# Copy the formatting of the row above (row 31, a normal task row) down into the
# newly inserted row so styles/number-formats/borders match the other task rows.
$ws.Range("A31:BZ31").Copy()
$ws.Range("A32:BZ32").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows(32).RowHeight = 30

$ws.Range("A32").Value2 = ""
$ws.Range("B32").Value = "add source codes"
$ws.Range("C32").Value = "Airat"
$ws.Range("D32").Value = 0
$ws.Range("E32").Formula = "=E30"
$ws.Range("F32").Formula = "=E32"

# --- Give the newly-inserted row the same day-count formula as every other task row ---
$ws.Range("H32").Formula = "=IF(OR(ISBLANK(task_start),ISBLANK(task_end)),"""",task_end-task_start+1)"

# --- Extend conditional-formatting ranges that used to stop at row 39 to row 40 ---
$ws.Range("D7:D39").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D7:D40"))
$ws.Range("I5:BL39").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I5:BL40"))
$ws.Range("I7:BL39").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I7:BL40"))
$ws.Range("I7:BL39").FormatConditions.Item(2).ModifyAppliesToRange($ws.Range("I7:BL40"))
$ws.Range("BM5:BS39").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("BM5:BS40"))
$ws.Range("BM7:BS39").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("BM7:BS40"))
$ws.Range("BM7:BS39").FormatConditions.Item(2).ModifyAppliesToRange($ws.Range("BM7:BS40"))
$ws.Range("BT5:BZ39").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("BT5:BZ40"))
$ws.Range("BT7:BZ39").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("BT7:BZ40"))
$ws.Range("BT7:BZ39").FormatConditions.Item(2).ModifyAppliesToRange($ws.Range("BT7:BZ40"))

# --- Update the frozen-pane view / active selection to match where the author left off ---
$ws.Application.Goto($ws.Range("A8"))
$ws.Range("E33").Select()
